# Fix "duplicate excel bug": Name/Is Student?/Likes Cats columns had been
# overwritten with duplicated values (both people showing "Hi " / "Jfjvk" /
# "Male" / "No"). Restore the real, distinct answers for each respondent and
# clear the stray duplicate that had been written into the "Favorite
# Subject" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Peter
$ws.Range("A3").Value = "Peter"
$ws.Range("B3").Value = "Yes"
$ws.Range("C3").Value = $null
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "None"
$ws.Range("G3").Value = "None"

# Row 4 - Meghan
$ws.Range("A4").Value = "Meghan"
$ws.Range("B4").Value = "No"
$ws.Range("C4").Value = $null
$ws.Range("E4").Value = "Yes"
$ws.Range("F4").Value = "None"
$ws.Range("G4").Value = "None"

# Column A was best-fit to its contents; re-fit now that "Peter"/"Meghan"
# replace the old duplicated "Hi " value so it stays sized to the text.
$ws.Columns.Item(1).AutoFit() | Out-Null
